{"js": "// DiaryService documentation: expand the \"Napl\u00f3z\u00e1s\" (logging) list item\n// into its own sub-list of CRUD operations, and add two more sub-items\n// after the existing \"Napl\u00f3z\u00e1s <-> SPA\" bullet.\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst items = paras.items;\nfor (const p of items) {\n  p.load(\"text\");\n}\nawait context.sync();\nfor (const p of items) {\n  p.load(\"listItemOrNullObject\");\n}\nawait context.sync();\nfor (const p of items) {\n  const li = p.listItemOrNullObject;\n  if (!li.isNullObject) {\n    li.load(\"level\");\n  }\n}\nawait context.sync();\n\n// Locate the two anchor paragraphs structurally:\n//   items[i]   : \"SPA ... Napl\u00f3z\u00e1s\" at list level 1 (w:ilvl=1)\n//   items[i+1] : \"Napl\u00f3z\u00e1s ... SPA\" at list level 1 (w:ilvl=1), immediately after\n// This pair is unique in the document, so we don't depend on fixed paragraph\n// index numbers.\nlet targetIdx = -1;\nfor (let i = 0; i < items.length - 1; i++) {\n  const t = items[i].text;\n  const tn = items[i + 1].text;\n  const li = items[i].listItemOrNullObject;\n  const lin = items[i + 1].listItemOrNullObject;\n  if (\n    !li.isNullObject &&\n    !lin.isNullObject &&\n    t.startsWith(\"SPA\") &&\n    t.trim().endsWith(\"Napl\u00f3z\u00e1s\") &&\n    tn.startsWith(\"Napl\u00f3z\u00e1s\") &&\n    tn.trim().endsWith(\"SPA\") &&\n    li.level === 1 &&\n    lin.level === 1\n  ) {\n    targetIdx = i;\n    break;\n  }\n}\n\nif (targetIdx === -1) {\n  throw new Error(\"anchor paragraphs not found\");\n}\n\n// 1) Promote the \"SPA -> Napl\u00f3z\u00e1s\" bullet from level 1 to level 2 (w:ilvl 1 -> 2)\nconst spaToNaplo = items[targetIdx];\nspaToNaplo.listItemOrNullObject.level = 2;\n\n// 2) Insert the four new sub-bullets (level 2 / w:ilvl=2) right after it\nconst newTexts = [\n  \"Napl\u00f3adatok lek\u00e9rdez\u00e9se\",\n  \"\u00daj bejegyz\u00e9s k\u00e9sz\u00edt\u00e9se\",\n  \"Bejegyz\u00e9s szerkeszt\u00e9se\",\n  \"Bejegyz\u00e9s t\u00f6rl\u00e9se\",\n];\nlet insertAfter = spaToNaplo;\nfor (const txt of newTexts) {\n  const np = insertAfter.insertParagraph(txt, Word.InsertLocation.after);\n  np.listItemOrNullObject.level = 2;\n  insertAfter = np;\n}\n\n// The \"Napl\u00f3z\u00e1s -> SPA\" paragraph (untouched, level 1) sits right after the\n// original promoted bullet (i.e. right after items[targetIdx+1] in the\n// original collection); it's now right after the four inserted bullets.\nconst naploToSpa = items[targetIdx + 1];\n\n// 3) Insert the two remaining new sub-bullets (level 2 / w:ilvl=2) after it\nconst moreTexts = [\"Napl\u00f3adatok visszaad\u00e1sa\", \"http v\u00e1lasz\u00fczenetek \"];\ninsertAfter = naploToSpa;\nfor (const txt of moreTexts) {\n  const np = insertAfter.insertParagraph(txt, Word.InsertLocation.after);\n  np.listItemOrNullObject.level = 2;\n  insertAfter = np;\n}\n\nawait context.sync();\n", "ps1": "# DiaryService documentation: expand the \"Napl\u00f3z\u00e1s\" (logging) list item\n# into its own sub-list of CRUD operations, and add two more sub-items\n# after the existing \"Napl\u00f3z\u00e1s <-> SPA\" bullet.\n$d = $word.ActiveDocument\n\n# Locate the two anchor paragraphs structurally:\n#   p1: \"SPA ... Napl\u00f3z\u00e1s\" at list level 2 (w:ilvl=1)\n#   p2 (immediately after p1): \"Napl\u00f3z\u00e1s ... SPA\" at list level 2 (w:ilvl=1)\n# This pair is unique in the document, so we don't depend on paragraph index\n# numbers that could shift between engines.\n$n = $d.Paragraphs.Count\n$targetIndex = -1\nfor ($i = 1; $i -lt $n; $i++) {\n    $p = $d.Paragraphs($i)\n    $pnext = $d.Paragraphs($i + 1)\n    $t = $p.Range.Text\n    $tn = $pnext.Range.Text\n    $cond1 = $t.StartsWith(\"SPA\") -and $t.Trim().EndsWith(\"Napl\u00f3z\u00e1s\")\n    $cond2 = $tn.StartsWith(\"Napl\u00f3z\u00e1s\") -and $tn.Trim().EndsWith(\"SPA\")\n    $cond3 = ($p.Range.ListFormat.ListLevelNumber -eq 2) -and ($pnext.Range.ListFormat.ListLevelNumber -eq 2)\n    if ($cond1 -and $cond2 -and $cond3) {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    Write-Output \"ERROR: anchor paragraphs not found\"\n} else {\n    # 1) Promote the \"SPA -> Napl\u00f3z\u00e1s\" bullet from level 2 to level 3 (w:ilvl 1 -> 2)\n    $p1 = $d.Paragraphs($targetIndex)\n    $p1.Range.ListFormat.ListLevelNumber = 3\n\n    # 2) Insert the four new sub-bullets (level 3 / w:ilvl=2) right after it\n    $newTexts = @(\"Napl\u00f3adatok lek\u00e9rdez\u00e9se\", \"\u00daj bejegyz\u00e9s k\u00e9sz\u00edt\u00e9se\", \"Bejegyz\u00e9s szerkeszt\u00e9se\", \"Bejegyz\u00e9s t\u00f6rl\u00e9se\")\n    $cur = $targetIndex\n    foreach ($txt in $newTexts) {\n        $d.Paragraphs($cur).Range.InsertParagraphAfter()\n        $cur = $cur + 1\n        $np = $d.Paragraphs($cur)\n        $np.Range.Text = $txt\n        $np.Range.ListFormat.ListLevelNumber = 3\n    }\n\n    # The \"Napl\u00f3z\u00e1s -> SPA\" paragraph (untouched, level 2) now sits right after\n    # the four inserted bullets.\n    $cur = $cur + 1\n    $pSpa = $d.Paragraphs($cur)\n\n    # 3) Insert the two remaining new sub-bullets (level 3 / w:ilvl=2) after it\n    $moreTexts = @(\"Napl\u00f3adatok visszaad\u00e1sa\", \"http v\u00e1lasz\u00fczenetek \")\n    foreach ($txt in $moreTexts) {\n        $pSpa.Range.InsertParagraphAfter()\n        $cur = $cur + 1\n        $np = $d.Paragraphs($cur)\n        $np.Range.Text = $txt\n        $np.Range.ListFormat.ListLevelNumber = 3\n        $pSpa = $np\n    }\n}\n"}
